# "display limit set to 10"
# The attendance table on Sheet1 previously only listed one 10-day cycle
# (rows 2-16). The display limit was raised so the same 10-day cycle
# (rows 11-16 hold the last 6 days of that cycle) now repeats three more
# times, extending the table from A1:E16 down to A1:E34.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the repeating 6-row block (rows 11-16) and paste it three more
# times directly below the existing data, preserving values, shared
# string references and cell styles exactly.
$ws.Range("A11:E16").Copy() | Out-Null
$ws.Range("A17:E22").PasteSpecial() | Out-Null
$ws.Range("A23:E28").PasteSpecial() | Out-Null
$ws.Range("A29:E34").PasteSpecial() | Out-Null

$excel.CutCopyMode = 0
